# Updated cryptos list with GitHub Actions.
# Refresh the Price (D) and Volume(1h) (E) columns for each coin row.
# Numeric-looking Price strings are entered with a leading apostrophe so
# Excel keeps them as literal text (matching the source data, e.g.
# "303.43" rather than being coerced into a number), then the style is
# reset to "Normal" so no stray number-format/quote-prefix style sticks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.238.18"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "1.603.10"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "'303.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("D7").Value = "'0.3765"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  +4.58%  "
$ws.Range("D9").Value = "'0.3631"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").Value = "'0.08120"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").Value = "'6.593"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "'7.414"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D16").Value = "'0.00001246"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "1.605.26"
$ws.Range("D18").Value = "'94.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.20%  "
$ws.Range("D19").Value = "'0.06930"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").Value = "'18.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").Value = "'6.529"
$ws.Range("D21").Style = "Normal"
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("D24").Value = "23.254.85"
$ws.Range("E24").Value = "  +1.04%  "
$ws.Range("D25").Value = "'3.034"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.16%  "
$ws.Range("D26").Value = "'2.371"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("D27").Value = "'21.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("D28").Value = "'149.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("D29").Value = "'5.260"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("D30").Value = "'134.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("E31").Value = "  +4.38%  "
$ws.Range("D32").Value = "'6.723"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("D33").Value = "1.782.49"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").Value = "'0.9614"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").Value = "'0.07498"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.75%  "
$ws.Range("D36").Value = "'0.02752"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("D37").Value = "'10.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("D38").Value = "'0.2537"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("D39").Value = "'6.115"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.84%  "
$ws.Range("D40").Value = "'0.08807"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D41").Value = "'1.395"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.23%  "
$ws.Range("D42").Value = "'0.7111"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("D43").Value = "'12.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "'15.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.25%  "
$ws.Range("D45").Value = "'0.6538"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("D46").Value = "'2.317"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").Value = "'4.015"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'132.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").Value = "'0.07949"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").Value = "'1.204"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.61%  "
